# Invoice 5 - add item/description/qty/price/amount header row and
# consolidate the "INVOICE 5" / "Uranka's Outdoors" banner cells into a
# single "Uranka's Outdoors - Invoice" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$g1 = $ws.Range("G1")
$h2 = $ws.Range("H2")
$e2 = $ws.Range("E2")

# Move/merge the banner text: "INVOICE 5" (G1) + "Uranka's Outdoors" (H2)
# become a single "Uranka's Outdoors - Invoice" cell at E2, keeping G1's
# look (big centered title font).
$e2.Value = "Uranka's Outdoors - Invoice"
$g1.Copy()
$e2.PasteSpecial(-4122)

$g1.Clear()
$h2.Clear()

# Give columns B:F the same 20-wide column width that column A already has,
# to make room for the new line-item table.
$ws.Range("B1:F1").ColumnWidth = 20

# New line-item header row.
$a11 = $ws.Range("A11")
$b11 = $ws.Range("B11")
$c11 = $ws.Range("C11")
$d11 = $ws.Range("D11")
$e11 = $ws.Range("E11")
$f11 = $ws.Range("F11")

$a11.Value = "Item"
$b11.Value = "Description"
$c11.Value = "Quantity"
$d11.Value = "Unit Price"
$e11.Value = "Amount"
$f11.Value = "Amount Due"

$hdr = $ws.Range("A11:F11")
$a2 = $ws.Range("A2")
$a2.Copy()
$hdr.PasteSpecial(-4122)
$hdr.Font.Size = 10
$hdr.HorizontalAlignment = 1
